$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the header cells that previously spanned multiple columns
# (Tkl group H1:L1, Challenges group M1:P1, Blocks group Q1:S1) so each
# column gets its own header label.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Fill in the visible header row (row 1) with a per-column label.
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# Row 2 (the original pandas sub-header) is kept but hidden, and a new
# blank hidden row 3 is inserted as a spacer.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true

# A handful of data rows were missing the "Tkl%" (O) column entirely;
# backfill them with 0 like their neighbours.
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("O19").Value = 0

# The totals row (20) is kept but hidden.
$ws.Rows.Item(20).Hidden = $true
